$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.734.61"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.76%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.166.18"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -4.73%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "589.83"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.59%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "134.58"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -5.84%  "
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.165.44"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -4.68%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.516"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.71%  "
$ws.Range("E10").Value = "  -6.85%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.23"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -5.87%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.451"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.97%  "
$ws.Range("E13").Value = "  -5.31%  "
$ws.Range("E14").Value = "  -1.05%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.686.61"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.80%  "
$ws.Range("E16").Value = "  -1.16%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.166.97"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -4.78%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "62.747.81"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.90%  "
$ws.Range("E19").Value = "  -4.88%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "458.93"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.90%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.98"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.84%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.691"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -6.57%  "
$ws.Range("E23").Value = "  -4.62%  "
$ws.Range("E24").Value = "  -4.44%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "82.79"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.60%  "
$ws.Range("E26").Value = "  -0.03%  "
$ws.Range("E27").Value = "  +0.00%  "
$ws.Range("E28").Value = "  -4.38%  "
$ws.Range("E29").Value = "  -7.20%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.76"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -5.99%  "
$ws.Range("E31").Value = "  -6.58%  "
$ws.Range("E32").Value = "  -6.34%  "
$ws.Range("E33").Value = "  -4.07%  "
$ws.Range("E34").Value = "  -7.02%  "
$ws.Range("E35").Value = "  -6.78%  "
$ws.Range("E36").Value = "  -4.95%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "51.22"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.19%  "
$ws.Range("E38").Value = "  -6.73%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "402.96"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -7.39%  "
$ws.Range("E41").Value = "  -3.87%  "
$ws.Range("E42").Value = "  -5.65%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.61"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -5.59%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.785.31"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -11.08%  "
$ws.Range("E45").Value = "  -6.79%  "
$ws.Range("E47").Value = "  -6.57%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "124.51"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.10%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "25.14"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.97%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "34.23"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -7.17%  "
